$wb = $excel.ActiveWorkbook

# Row 33 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 183.6
$ws.Range("I33").Value = 180.57143
$ws.Range("J33").Value = 190.66667
$ws.Range("K33").Value = 180.57143
$ws.Range("L33").Value = 190.66667
$ws.Range("M33").Value = 48.42857000000001
$ws.Range("N33").Value = -648.6666700000001

# Row 64 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5228.263
$ws.Range("J64").Value = 3959.5454
$ws.Range("L64").Value = 3959.5454
$ws.Range("N64").Value = -4455.5454

# Row 67 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5228.263
$ws.Range("J67").Value = 3959.5454
$ws.Range("L67").Value = 3959.5454
$ws.Range("N67").Value = -5675.5454

# Row 74 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3920513.8
$ws.Range("I74").Value = 4246806.5
$ws.Range("K74").Value = 4246806.5
$ws.Range("M74").Value = -4245870.5

# Row 77 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3920513.8
$ws.Range("I77").Value = 4246806.5
$ws.Range("K77").Value = 21234032.5
$ws.Range("M77").Value = -21229352.5

# Row 100 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1565.5
$ws.Range("I100").Value = 1563.125
$ws.Range("J100").Value = 1575
$ws.Range("K100").Value = 1563.125
$ws.Range("L100").Value = 1575
$ws.Range("M100").Value = -1022.125
$ws.Range("N100").Value = -2657

# Row 125 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 8188.6665
$ws.Range("I125").Value = 3170
$ws.Range("J125").Value = 13207.333
$ws.Range("K125").Value = 28530
$ws.Range("L125").Value = 118865.997
$ws.Range("M125").Value = -26070
$ws.Range("N125").Value = -123785.997

# Row 129 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 950644.1
$ws.Range("I129").Value = 446.3
$ws.Range("J129").Value = 1278298.5
$ws.Range("K129").Value = 1338.9
$ws.Range("L129").Value = 3834895.5
$ws.Range("M129").Value = 3661.1
$ws.Range("N129").Value = -3844895.5

# Row 139 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 53975
$ws.Range("J139").Value = 53975
$ws.Range("L139").Value = 53975
$ws.Range("N139").Value = -64255

# Row 2 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3020.348
$ws.Range("I2").Value = 2522.2144
$ws.Range("J2").Value = 3795.2222
$ws.Range("K2").Value = 2522.2144
$ws.Range("L2").Value = 3795.2222
$ws.Range("M2").Value = -2409.2144
$ws.Range("N2").Value = -4021.2222

# Row 61 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 777.86206
$ws.Range("I61").Value = 736.3077
$ws.Range("K61").Value = 736.3077
$ws.Range("M61").Value = -524.3077

# Row 102 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4534.9
$ws.Range("I102").Value = 4154.364
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 4154.364
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -2532.364
$ws.Range("N102").Value = -8244

# Row 116 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3020.348
$ws.Range("I116").Value = 2522.2144
$ws.Range("J116").Value = 3795.2222
$ws.Range("K116").Value = 2522.2144
$ws.Range("L116").Value = 3795.2222
$ws.Range("M116").Value = -228.2143999999998
$ws.Range("N116").Value = -8383.2222

# Row 122 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 952
$ws.Range("I122").Value = 945.1429000000001
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2835.4287
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -385.4287000000004
$ws.Range("N122").Value = -7900

# Row 136 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 777.86206
$ws.Range("I136").Value = 736.3077
$ws.Range("K136").Value = 2208.9231
$ws.Range("M136").Value = 341.0769

# Row 3 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3020.348
$ws.Range("I3").Value = 2522.2144
$ws.Range("J3").Value = 3795.2222
$ws.Range("K3").Value = 2522.2144
$ws.Range("L3").Value = 3795.2222
$ws.Range("M3").Value = -2408.2144
$ws.Range("N3").Value = -4023.2222

# Row 22 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 289.5
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 289.5
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 289.5
$ws.Range("N22").Value = -635.5
$ws.Range("M22").ClearContents()

# Row 94 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1484.1428
$ws.Range("I94").Value = 1155.6
$ws.Range("J94").Value = 1666.6666
$ws.Range("K94").Value = 1155.6
$ws.Range("L94").Value = 1666.6666
$ws.Range("M94").Value = -704.5999999999999
$ws.Range("N94").Value = -2568.6666

# Row 105 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5438.2915
$ws.Range("I105").Value = 5436.4116
$ws.Range("J105").Value = 5442.857
$ws.Range("K105").Value = 5436.4116
$ws.Range("L105").Value = 5442.857
$ws.Range("M105").Value = -3689.4116
$ws.Range("N105").Value = -8936.857

# Row 31 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2443.861
$ws.Range("I31").Value = 2411.147
$ws.Range("K31").Value = 2411.147
$ws.Range("M31").Value = -2116.147

# Row 34 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2443.861
$ws.Range("I34").Value = 2411.147
$ws.Range("K34").Value = 2411.147
$ws.Range("M34").Value = -2209.147

# Row 58 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8411.6
$ws.Range("I58").Value = 1655.2858
$ws.Range("J58").Value = 103000
$ws.Range("K58").Value = 1655.2858
$ws.Range("L58").Value = 103000
$ws.Range("M58").Value = -1452.2858
$ws.Range("N58").Value = -103406

# Row 62 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 9400
$ws.Range("I62").Value = 9480
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 9480
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = -8856
$ws.Range("N62").Value = -10248

# Row 65 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 9400
$ws.Range("I65").Value = 9480
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 47400
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = -44280
$ws.Range("N65").Value = -51240

# Row 136 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8411.6
$ws.Range("I136").Value = 1655.2858
$ws.Range("J136").Value = 103000
$ws.Range("K136").Value = 4965.857400000001
$ws.Range("L136").Value = 309000
$ws.Range("M136").Value = -2415.857400000001
$ws.Range("N136").Value = -314100

# Row 80 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6607.6924
$ws.Range("J80").Value = 7490
$ws.Range("L80").Value = 22470
$ws.Range("N80").Value = -24342

# Row 83 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 6607.6924
$ws.Range("J83").Value = 7490
$ws.Range("L83").Value = 67410
$ws.Range("N83").Value = -76770

# Row 22 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

# Row 70 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# Row 73 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# Row 132 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3343.6428
$ws.Range("I132").Value = 3283.7646
$ws.Range("J132").Value = 3436.182
$ws.Range("K132").Value = 9851.293799999999
$ws.Range("L132").Value = 10308.546
$ws.Range("M132").Value = -7321.293799999999
$ws.Range("N132").Value = -15368.546

# Row 138 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 55812.25
$ws.Range("I138").Value = 35000
$ws.Range("J138").Value = 62749.668
$ws.Range("K138").Value = 35000
$ws.Range("L138").Value = 62749.668
$ws.Range("M138").Value = -29860
$ws.Range("N138").Value = -73029.66800000001

# Row 46 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1425.4546
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 1000
$ws.Range("N46").Value = -1376

# Row 100 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 80925
$ws.Range("I100").Value = 112685
$ws.Range("K100").Value = 112685
$ws.Range("M100").Value = -112144

# Row 132 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2008.6666
$ws.Range("I132").Value = 1410.5333
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 4231.5999
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -1701.5999
$ws.Range("N132").Value = -20058.0005

# Row 138 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 39175
$ws.Range("J138").Value = 39175
$ws.Range("L138").Value = 39175
$ws.Range("N138").Value = -49455

# Row 96 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3826.6667
$ws.Range("I96").Value = 3822.2222
$ws.Range("J96").Value = 3833.3333
$ws.Range("K96").Value = 3822.2222
$ws.Range("L96").Value = 3833.3333
$ws.Range("M96").Value = -2449.2222
$ws.Range("N96").Value = -6579.3333

# Row 100 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1199.6666
$ws.Range("I100").Value = 338.2
$ws.Range("J100").Value = 2276.5
$ws.Range("K100").Value = 676.4
$ws.Range("L100").Value = 4553
$ws.Range("M100").Value = -135.4
$ws.Range("N100").Value = -5635

# Row 133 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 66948.5
$ws.Range("J133").Value = 66948.5
$ws.Range("L133").Value = 66948.5
$ws.Range("N133").Value = -77068.5
